$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" (column D) values look like plain numbers (e.g. "586.55")
# and Excel's COM Value setter would silently coerce those into numeric
# cells, losing the original text-cell representation used throughout
# this sheet. Force those specific cells to Text format first so the
# assigned strings are stored verbatim as text, matching the source data.
# (Cells whose price already contains multiple '.' separators, e.g.
# "70.183.83", are never auto-converted, so they're left alone.)
$numericLookingPriceCells = @(
  "D5","D6","D7","D10","D11","D12","D13","D14","D17","D19",
  "D22","D23","D24","D25","D26","D27","D28","D29","D30","D31",
  "D32","D33","D35","D38","D39","D41","D43","D44","D45","D46",
  "D48","D49","D50"
)
foreach ($addr in $numericLookingPriceCells) {
  $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "69.965.60"
$ws.Range("E2").Value = "  +4.02%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.594.87"
$ws.Range("E3").Value = "  +4.12%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - BNB
$ws.Range("D5").Value = "586.55"
$ws.Range("E5").Value = "  +3.20%  "

# Row 6 - Solana
$ws.Range("D6").Value = "189.55"
$ws.Range("E6").Value = "  +2.67%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.641"
$ws.Range("E7").Value = "  +1.18%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.591.56"
$ws.Range("E8").Value = "  +4.19%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  +0.00%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.175"
$ws.Range("E10").Value = "  -0.86%  "

# Row 11 - Cardano
$ws.Range("D11").Value = "0.657"
$ws.Range("E11").Value = "  +1.75%  "

# Row 12 - Avalanche
$ws.Range("D12").Value = "57.65"
$ws.Range("E12").Value = "  +3.68%  "

# Row 13 - ShibaInu
$ws.Range("D13").Value = "0.0000288"
$ws.Range("E13").Value = "  +2.35%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "9.74"
$ws.Range("E14").Value = "  +3.99%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.182.23"
$ws.Range("E15").Value = "  +4.32%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.604.40"
$ws.Range("E16").Value = "  +4.42%  "

# Row 17 - Chainlink
$ws.Range("D17").Value = "19.32"
$ws.Range("E17").Value = "  +4.15%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "70.000.07"
$ws.Range("E18").Value = "  +3.99%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "12.43"
$ws.Range("E19").Value = "  +3.40%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  +0.31%  "

# Row 21 - Polygon
$ws.Range("E21").Value = "  +3.33%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "488.43"
$ws.Range("E22").Value = "  +0.93%  "

# Row 23 - InternetComputer(DFINITY)
$ws.Range("D23").Value = "17.44"
$ws.Range("E23").Value = "  +15.63%  "

# Row 24 - Toncoin
$ws.Range("D24").Value = "5.37"
$ws.Range("E24").Value = "  +8.44%  "

# Row 25 - PancakeSwap
$ws.Range("D25").Value = "4.43"
$ws.Range("E25").Value = "  +6.08%  "

# Row 26 - Litecoin
$ws.Range("D26").Value = "90.34"
$ws.Range("E26").Value = "  +0.57%  "

# Row 27 - ImmutableX
$ws.Range("D27").Value = "3.09"
$ws.Range("E27").Value = "  +4.34%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "11.02"
$ws.Range("E28").Value = "  +1.02%  "

# Row 29 - Filecoin
$ws.Range("D29").Value = "9.38"
$ws.Range("E29").Value = "  +4.92%  "

# Row 30 - EthereumClassic
$ws.Range("D30").Value = "32.21"
$ws.Range("E30").Value = "  +1.95%  "

# Row 31 - NEARProtocol
$ws.Range("D31").Value = "7.47"
$ws.Range("E31").Value = "  +7.07%  "

# Row 32 - Bittensor
$ws.Range("D32").Value = "624.50"
$ws.Range("E32").Value = "  +4.18%  "

# Row 33 - Cosmos
$ws.Range("D33").Value = "12.19"
$ws.Range("E33").Value = "  +4.98%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  +6.32%  "

# Row 35 - OKB
$ws.Range("D35").Value = "65.02"
$ws.Range("E35").Value = "  +3.18%  "

# Row 36 - PEPE
$ws.Range("D36").Value = "0.0₃0810"
$ws.Range("E36").Value = "  +3.70%  "

# Row 37 - Dai
$ws.Range("E37").Value = "  +0.06%  "

# Row 38 - TheGraph
$ws.Range("D38").Value = "0.402"
$ws.Range("E38").Value = "  +3.61%  "

# Row 39 - InjectiveProtocol
$ws.Range("D39").Value = "37.78"
$ws.Range("E39").Value = "  +3.44%  "

# Row 40 - Kaspa
$ws.Range("E40").Value = "  -1.32%  "

# Row 41 - Stacks
$ws.Range("D41").Value = "3.61"
$ws.Range("E41").Value = "  -1.00%  "

# Row 42 - Maker
$ws.Range("D42").Value = "3.299.36"
$ws.Range("E42").Value = "  +5.02%  "

# Row 43 - ThetaToken
$ws.Range("D43").Value = "3.06"
$ws.Range("E43").Value = "  +4.54%  "

# Row 44 - VeChain
$ws.Range("D44").Value = "0.0442"
$ws.Range("E44").Value = "  +3.79%  "

# Row 45 - Fetch.AI
$ws.Range("D45").Value = "2.65"
$ws.Range("E45").Value = "  +2.65%  "

# Row 46 - ApeXProtocol
$ws.Range("D46").Value = "3.33"
$ws.Range("E46").Value = "  +2.62%  "

# Row 47 - Stellar
$ws.Range("E47").Value = "  +1.26%  "

# Row 48 - THORChain
$ws.Range("D48").Value = "9.08"
$ws.Range("E48").Value = "  +3.81%  "

# Row 49 - was dogwifhat, now LidoDAOToken (rows 49/50 swapped content)
$ws.Range("B49").Value = "LidoDAOToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D49").Value = "3.32"
$ws.Range("E49").Value = "  +5.89%  "

# Row 50 - was LidoDAOToken, now dogwifhat
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "2.68"
$ws.Range("E50").Value = "  -4.71%  "

# Row 51 - FirstDigitalUSD
$ws.Range("E51").Value = "  -0.01%  "
